# Regenerate orders with updated distance/size labels.
# The experiment's distance and size condition codes changed:
#   D80 -> D86, D51 -> D55, D64 -> D69, S30 -> S31
# These codes appear embedded inside many cell values (condition names,
# filenames, distance labels, size labels) throughout the sheet, so we
# perform a global substring Find & Replace across every used cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.UsedRange

$cells.Replace("D80", "D86")
$cells.Replace("D51", "D55")
$cells.Replace("D64", "D69")
$cells.Replace("S30", "S31")
